$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2806.6667
$ws.Range("I40").Value = 1490
$ws.Range("J40").Value = 5440
$ws.Range("K40").Value = 1490
$ws.Range("L40").Value = 5440
$ws.Range("M40").Value = -1315
$ws.Range("N40").Value = -5790

$ws.Range("H112").Value = 1766.6666
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 2000
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -8216

$ws.Range("H118").Value = 296
$ws.Range("I118").Value = 225
$ws.Range("K118").Value = 675
$ws.Range("M118").Value = 982

$ws.Range("H138").Value = 2627

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 137.11111
$ws.Range("I5").Value = 120.42857
$ws.Range("J5").Value = 195.5
$ws.Range("K5").Value = 120.42857
$ws.Range("L5").Value = 195.5
$ws.Range("M5").Value = -8.428569999999993
$ws.Range("N5").Value = -419.5

$ws.Range("H22").Value = 1666.6666
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 4000
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 4000
$ws.Range("M22").Value = -201
$ws.Range("N22").Value = -4598

$ws.Range("H61").Value = 4199.5
$ws.Range("I61").Value = 4199.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4199.5
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -3987.5
$ws.Range("M61").ClearContents()

$ws.Range("H136").Value = 4199.5
$ws.Range("I136").Value = 4199.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12598.5
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -10048.5
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 137.11111
$ws.Range("I4").Value = 120.42857
$ws.Range("J4").Value = 195.5
$ws.Range("K4").Value = 120.42857
$ws.Range("L4").Value = 195.5
$ws.Range("M4").Value = -5.428569999999993
$ws.Range("N4").Value = -425.5

$ws.Range("H22").Value = 212.75
$ws.Range("I22").Value = 250.33333
$ws.Range("K22").Value = 250.33333
$ws.Range("M22").Value = -77.33332999999999

$ws.Range("H134").Value = 15000
$ws.Range("I134").Value = 15000
$ws.Range("K134").Value = 45000
$ws.Range("M134").Value = -42465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("L4").ClearContents()

$ws.Range("H7").Value = 410.35294
$ws.Range("I7").Value = 414.5
$ws.Range("K7").Value = 414.5
$ws.Range("M7").Value = -301.5

$ws.Range("H22").Value = 374.25
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = 100

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("N58").Value = 0
$ws.Range("L58").ClearContents()
$ws.Range("M58").ClearContents()

$ws.Range("H99").Value = 7000
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 7000
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = 7000
$ws.Range("N99").Value = -9996
$ws.Range("L99").ClearContents()

$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = 21000
$ws.Range("N126").Value = -25940
$ws.Range("L126").ClearContents()

$ws.Range("H132").Value = 2026.5
$ws.Range("I132").Value = 289.75
$ws.Range("K132").Value = 869.25
$ws.Range("M132").Value = 1660.75

$ws.Range("H135").Value = 49832.832
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("N136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 265.625
$ws.Range("J2").Value = 541.6667
$ws.Range("L2").Value = 3250.0002
$ws.Range("N2").Value = -3476.0002

$ws.Range("H4").Value = 47625.668
$ws.Range("I4").Value = 75014.47
$ws.Range("J4").Value = 1977.6666
$ws.Range("K4").Value = 225043.41
$ws.Range("L4").Value = 5932.9998
$ws.Range("M4").Value = -224931.41
$ws.Range("N4").Value = -6156.9998

$ws.Range("H13").Value = 90
$ws.Range("I13").Value = 80
$ws.Range("K13").Value = 240
$ws.Range("M13").Value = -72

$ws.Range("H15").Value = 50
$ws.Range("J15").Value = 50
$ws.Range("L15").Value = 150
$ws.Range("N15").Value = -430

$ws.Range("H16").Value = 422.25
$ws.Range("I16").Value = 396.33334
$ws.Range("K16").Value = 1189.00002
$ws.Range("M16").Value = -1016.00002

$ws.Range("H20").Value = 291.66666
$ws.Range("I20").Value = 237.5
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 712.5
$ws.Range("L20").Value = 1200
$ws.Range("M20").Value = -485.5
$ws.Range("N20").Value = -1654

$ws.Range("H21").Value = 300
$ws.Range("I21").Value = 300
$ws.Range("K21").Value = 900
$ws.Range("M21").Value = -727

$ws.Range("H80").Value = 7985.75
$ws.Range("J80").Value = 7814.6665
$ws.Range("L80").Value = 23443.9995
$ws.Range("N80").Value = -25315.9995

$ws.Range("H83").Value = 7985.75
$ws.Range("J83").Value = 7814.6665
$ws.Range("L83").Value = 70331.9985
$ws.Range("N83").Value = -79691.9985

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 21000
$ws.Range("I18").Value = 3000
$ws.Range("J18").Value = 39000
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 39000
$ws.Range("M18").Value = -2707
$ws.Range("N18").Value = -39586

$ws.Range("H102").Value = 5457.375
$ws.Range("I102").Value = 5524.143
$ws.Range("J102").Value = 4990
$ws.Range("K102").Value = 5524.143
$ws.Range("L102").Value = 4990
$ws.Range("M102").Value = -3902.143
$ws.Range("N102").Value = -8234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 59245
$ws.Range("J34").Value = 59245
$ws.Range("L34").Value = 59245
$ws.Range("N34").Value = -59589

$ws.Range("H40").Value = 5381
$ws.Range("I40").Value = 4800
$ws.Range("J40").Value = 6252.5
$ws.Range("K40").Value = 4800
$ws.Range("L40").Value = 6252.5
$ws.Range("M40").Value = -4664
$ws.Range("N40").Value = -6524.5

$ws.Range("H46").Value = 2568.1365
$ws.Range("I46").Value = 2275.0833
$ws.Range("J46").Value = 2919.8
$ws.Range("K46").Value = 2275.0833
$ws.Range("L46").Value = 2919.8
$ws.Range("M46").Value = -2087.0833
$ws.Range("N46").Value = -3295.8

$ws.Range("H132").Value = 2799.75
$ws.Range("I132").Value = 2799.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8399.25
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -5869.25
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 93482.164
$ws.Range("I126").Value = 90450.5
$ws.Range("J126").Value = 94998
$ws.Range("K126").Value = 271351.5
$ws.Range("L126").Value = 284994
$ws.Range("M126").Value = -268881.5
$ws.Range("N126").Value = -289934

$ws.Range("H132").Value = 4004
$ws.Range("I132").Value = 4004
$ws.Range("K132").Value = 12012
$ws.Range("M132").Value = -9482
